$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) stays text, matching the original formatting,
# so numeric-looking strings like "242.71" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.896.85"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "1.888.69"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "0.7668"
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("D6").Value = "242.71"
$ws.Range("E6").Value = "  -0.58%  "

$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "0.3134"
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "25.70"
$ws.Range("E9").Value = "  +1.59%  "

$ws.Range("D10").Value = "0.07171"
$ws.Range("E10").Value = "  -2.93%  "

$ws.Range("D11").Value = "0.08506"
$ws.Range("E11").Value = "  +4.46%  "

$ws.Range("D12").Value = "0.7633"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.373"
$ws.Range("E13").Value = "  -1.87%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.851.66"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("D15").Value = "93.91"
$ws.Range("E15").Value = "  +1.51%  "

$ws.Range("D16").Value = "6.146"
$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").Value = "29.967.26"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Value = "13.79"
$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("D19").Value = "244.31"
$ws.Range("E19").Value = "  -0.41%  "

$ws.Range("D20").Value = "0.000007807"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.125.54"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "8.066"
$ws.Range("E23").Value = "  -0.80%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").Value = "0.1628"
$ws.Range("E25").Value = "  +3.84%  "

$ws.Range("D26").Value = "9.400"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").Value = "162.29"
$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").Value = "18.77"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").Value = "1.494"
$ws.Range("E30").Value = "  +2.60%  "

$ws.Range("D31").Value = "1.539"
$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").Value = "4.489"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("D33").Value = "4.097"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").Value = "0.05425"
$ws.Range("E34").Value = "  -3.12%  "

$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("D36").Value = "0.7450"
$ws.Range("E36").Value = "  -1.66%  "

$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.28%  "

$ws.Range("D38").Value = "2.691"
$ws.Range("E38").Value = "  +1.63%  "

$ws.Range("D39").Value = "0.01949"
$ws.Range("E39").Value = "  +0.63%  "

$ws.Range("D40").Value = "2.782"
$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("D41").Value = "0.4468"
$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").Value = "1.100.05"
$ws.Range("E42").Value = "  -4.24%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "6.073"
$ws.Range("E43").Value = "  +1.83%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "72.90"
$ws.Range("E44").Value = "  -1.62%  "

$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").Value = "103.03"
$ws.Range("E47").Value = "  +1.01%  "

$ws.Range("D48").Value = "1.873"
$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("D49").Value = "7.665"
$ws.Range("E49").Value = "  +2.00%  "

$ws.Range("D50").Value = "2.995"
$ws.Range("E50").Value = "  -4.01%  "

$ws.Range("D51").Value = "2.018.39"
$ws.Range("E51").Value = "  -0.43%  "
